$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("A4").Value = 130789471
$ws.Range("B4").Value = 79243
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 490498
$ws.Range("R4").Value = 6763669
$ws.Range("Z4").Value = "14:18"
$ws.Range("AB4").Value = "14:18"
$ws.Range("AF4").Value = ""  # becomes a blank cell (new in target, matches empty inlineStr)
$ws.Range("AW4").Value = "Bo karlstens"
$ws.Range("AX4").Value = "Bo karlstens, Håkan Thenander"

# --- Row 5 ---
$ws.Range("A5").Value = 130754287
$ws.Range("B5").Value = 57881
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 490501
$ws.Range("R5").Value = 6763773
$ws.Range("Z5").Value = "11:43"
$ws.Range("AB5").Value = "11:43"
$ws.Range("AF5").ClearContents()
$ws.Range("AW5").Value = "Håkan Thenander"
$ws.Range("AX5").Value = "Håkan Thenander, Bo karlstens"

# --- Row 12 ---
$ws.Range("A12").Value = 130758082
$ws.Range("B12").Value = 79243
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value = 490186
$ws.Range("R12").Value = 6763602

# --- Row 13 ---
$ws.Range("A13").Value = 130757655
$ws.Range("B13").Value = 8451
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 106545
$ws.Range("F13").Value = "Mindre märgborre"
$ws.Range("G13").Value = "Tomicus minor"
$ws.Range("H13").Value = "(Hartig, 1834)"
$ws.Range("M13").Value = "äldre gnagspår"
$ws.Range("Q13").Value = 490214
$ws.Range("R13").Value = 6763590

# --- Row 35 ---
$ws.Range("A35").Value = 130757247
$ws.Range("B35").Value = 5197
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 105930
$ws.Range("F35").Value = "Vågbandad barkbock"
$ws.Range("G35").Value = "Semanotus undatus"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("M35").Value = "äldre gnagspår"
$ws.Range("Q35").Value = 490467
$ws.Range("R35").Value = 6763573
$ws.Range("Z35").Value = "15:01"
$ws.Range("AB35").Value = "15:01"
$ws.Range("AF35").ClearContents()
$ws.Range("AW35").Value = "Håkan Thenander"
$ws.Range("AX35").Value = "Håkan Thenander, Bo karlstens"

# --- Row 36 ---
$ws.Range("A36").Value = 130789468
$ws.Range("Q36").Value = 490321
$ws.Range("R36").Value = 6763593
$ws.Range("Z36").Value = "15:14"
$ws.Range("AB36").Value = "15:14"
$ws.Range("AC36").ClearContents()
$ws.Range("AF36").Value = ""  # becomes a blank cell (new in target, matches empty inlineStr)
$ws.Range("AW36").Value = "Bo karlstens"
$ws.Range("AX36").Value = "Bo karlstens, Håkan Thenander"

# --- Row 37 ---
$ws.Range("A37").Value = 130754851
$ws.Range("Q37").Value = 490449
$ws.Range("R37").Value = 6763949
$ws.Range("Z37").Value = "11:43"
$ws.Range("AB37").Value = "11:43"
$ws.Range("AC37").Value = "1 bild. På tallstam"

# --- Row 38 ---
$ws.Range("A38").Value = 130757412
$ws.Range("B38").Value = 79243
$ws.Range("D38").Value = "NT"
$ws.Range("E38").Value = 6425
$ws.Range("F38").Value = "Garnlav"
$ws.Range("G38").Value = "Alectoria sarmentosa"
$ws.Range("H38").Value = "(Ach.) Ach."
$ws.Range("M38").ClearContents()
$ws.Range("Q38").Value = 490381
$ws.Range("R38").Value = 6763583
